$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M7").Value = 2
$ws.Range("M8").Value = 2
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 2
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 2
$ws.Range("M11").Value = 2
$ws.Range("N11").Value = 2
$ws.Range("M12").Value = 10
$ws.Range("M13").Value = 4
$ws.Range("N14").Value = 4
$ws.Range("N16").Value = 4
